$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (these "price" cells use dotted
    # grouping that looks numeric, e.g. "23.534.37") and strip the
    # quote-prefix formatting that Excel applies for text-forced numbers so
    # the cell style is left untouched.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "23.534.37"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.651.04"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.83%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.45%  "

# Row 6 - BNB
Set-TextValue "D6" "300.44"
$ws.Range("E6").Value = "  -1.22%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3786"
$ws.Range("E7").Value = "  +0.21%  "

# Row 8 - was Cardano, now OKB
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D8" "50.80"
$ws.Range("E8").Value = "  -1.68%  "

# Row 9 - was OKB, now Cardano
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.3574"
$ws.Range("E9").Value = "  -0.76%  "

# Row 10 - Polygon
$ws.Range("E10").Value = "  -0.74%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.08115"
$ws.Range("E11").Value = "  -1.13%  "

# Row 12 - BinanceUSD
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.84%  "

# Row 13 - Solana
Set-TextValue "D13" "22.15"
$ws.Range("E13").Value = "  -1.40%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.429"
$ws.Range("E14").Value = "  -1.74%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.437"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.94%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.654.20"
$ws.Range("E17").Value = "  +0.67%  "

# Row 18 - Litecoin
Set-TextValue "D18" "97.23"
$ws.Range("E18").Value = "  +0.86%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06990"
$ws.Range("E19").Value = "  +0.30%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.793"
$ws.Range("E20").Value = "  +0.63%  "

# Row 21 - Avalanche
Set-TextValue "D21" "17.49"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22 - Dai
Set-TextValue "D22" "0.9998"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23 - Cosmos
Set-TextValue "D23" "12.64"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24 - WrappedBTC
Set-TextValue "D24" "23.562.19"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.486"
$ws.Range("E25").Value = "  -1.17%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.933"
$ws.Range("E26").Value = "  -6.51%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "21.03"
$ws.Range("E27").Value = "  -1.15%  "

# Row 28 - Monero
Set-TextValue "D28" "152.50"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - HuobiToken
Set-TextValue "D29" "5.233"
$ws.Range("E29").Value = "  +0.99%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "133.34"
$ws.Range("E30").Value = "  -0.78%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-TextValue "D31" "1.831.91"
$ws.Range("E31").Value = "  +0.45%  "

# Row 32 - Filecoin
Set-TextValue "D32" "7.021"
$ws.Range("E32").Value = "  +3.75%  "

# Row 33 - WEMIXTOKEN
$ws.Range("E33").Value = "  +5.55%  "

# Row 34 - FraxShare
Set-TextValue "D34" "11.90"
$ws.Range("E34").Value = "  +3.05%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "1.036"
$ws.Range("E35").Value = "  -5.25%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.02747"
$ws.Range("E36").Value = "  -1.11%  "

# Row 37 - Stellar
Set-TextValue "D37" "0.08711"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38 - InternetComputer(DFINITY)
Set-TextValue "D38" "6.012"
$ws.Range("E38").Value = "  -0.47%  "

# Row 39 - Algorand
Set-TextValue "D39" "0.2455"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40 - Aptos
Set-TextValue "D40" "13.22"
$ws.Range("E40").Value = "  +4.05%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.06898"
$ws.Range("E41").Value = "  -1.56%  "

# Row 42 - TheSandbox
Set-TextValue "D42" "0.6939"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43 - TrustWalletToken (E43 unchanged)
Set-TextValue "D43" "1.324"

# Row 44 - EnergySwap
Set-TextValue "D44" "15.71"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45 - Decentraland
Set-TextValue "D45" "0.6458"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46 - Frax
$ws.Range("E46").Value = "  +0.20%  "

# Row 47 - NEARProtocol
Set-TextValue "D47" "2.275"
$ws.Range("E47").Value = "  -2.52%  "

# Row 48 - PancakeSwap
Set-TextValue "D48" "3.933"
$ws.Range("E48").Value = "  -0.68%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.07821"
$ws.Range("E49").Value = "  -1.96%  "

# Row 50 - Quant
Set-TextValue "D50" "128.42"
$ws.Range("E50").Value = "  +0.46%  "

# Row 51 - Flow
$ws.Range("E51").Value = "  -0.54%  "
